# Generate Report for Handback
# This script fills in the handback columns (Latest Target File, Latest
# Handback File, Latest Handback DateTime) on the zh-cn and de-de report
# sheets, adds the matching "a.md" hyperlinks, refreshes the roll-up
# status on the Overview sheet, and widens the columns that now hold the
# longer status/file-name text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$statusText = "Handed back: in sync with en-US"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34b14e42d6562570b96a5e16776a0af7fa6475ce/e2e/a.md"
$bMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34b14e42d6562570b96a5e16776a0af7fa6475ce/e2e/b.md"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-locale status column for both rows.
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Range("E1").ColumnWidth = 29.12
$overview.Range("F1").ColumnWidth = 29.12

# ---------------------------------------------------------------------
# zh-cn sheet: record the handback (target file, handback file, handback
# datetime) for both rows, and link the new "Latest Target File" cell.
# ---------------------------------------------------------------------
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-01 06:41:42"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-01 06:41:42"

$zhcn.Range("C1").ColumnWidth = 29.12
$zhcn.Range("J1").ColumnWidth = 39.12

# ---------------------------------------------------------------------
# de-de sheet: record the handback (target file, handback file, handback
# datetime) for both rows, and link the new "Latest Target File" cell.
# ---------------------------------------------------------------------
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-09-01 06:41:49"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-09-01 06:41:49"

$dede.Range("C1").ColumnWidth = 29.12
$dede.Range("J1").ColumnWidth = 39.12

# ---------------------------------------------------------------------
# Hyperlinks: rebuild in display order (A2, I2, A3, I3) on each report
# sheet so the new "a.md" links on column I (Latest Target File) land
# next to the existing ones, matching the order Excel assigns rIds in.
# ---------------------------------------------------------------------
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $aMdUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aMdUrl, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $bMdUrl, "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aMdUrl, "", "", "a.md")

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $aMdUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I2"), $aMdUrl, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $bMdUrl, "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $aMdUrl, "", "", "a.md")
